$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A5").Value = 20240523
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 5
